$d = $word.ActiveDocument

# Locate the two target paragraphs by their (language-agnostic, ASCII-safe)
# text fragments instead of hard-coded indices, so the script keeps working
# even if paragraph numbering shifts a little.
$paraAnalog = $null
$paraThreshold = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*wczytanie ustawie*") {
        $paraAnalog = $p
    }
    if ($t -like "*Progowanie warto*cyfrowych*") {
        $paraThreshold = $p
    }
}

# The stray "_GoBack" bookmark (Word's "last edit position" marker) currently
# sits around the "Wymiana plikow konfiguracyjnych na INI" paragraph; it needs
# to move to wrap the "wczytanie ustawien ..." paragraph, matching a fresh
# edit having been made there.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Apply strikethrough formatting to both newly-finished/obsolete items.
$paraAnalog.Range.Font.StrikeThrough = 1
$paraThreshold.Range.Font.StrikeThrough = 1

# Re-insert the "_GoBack" bookmark spanning the "wczytanie ustawien ..."
# paragraph (including its paragraph mark), reflecting it being the most
# recently edited location.
$d.Bookmarks.Add("_GoBack", $paraAnalog.Range)
